$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: ID, First Name, Last Name (inherits the highlighted header style
# from the existing A1 cell)
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "First Name"
$ws.Range("C1").Value = "Last Name"

# First employee row
$c = $ws.Range("A2")
$c.NumberFormat = "@"
$c.Value = "0312"
$c.Style = "Normal"

$ws.Range("B3").Value = "A8DCo 4Ys"
$ws.Range("C4").Value = "010Z"

$ws.Range("F12").Select() | Out-Null
